# Weekly data refresh: insert a new most-recent price observation as a new
# row 164 in the "Arveja Verde" price history sheet, pushing the existing
# rows 164-178 down to 165-179 (dimension grows from R178 to R179).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 164, shifting rows 164:178 down to 165:179.
$ws.Rows(164).Insert()

# Populate the new row 164 with the latest weekly record.
$ws.Range("A164").Value = 4
$ws.Range("B164").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C164").Value = "Los Lagos"
$ws.Range("D164").Value = 45223
$ws.Range("E164").Value = 10
$ws.Range("F164").Value = 100112022
$ws.Range("G164").Value = "Arveja Verde"
$ws.Range("H164").Value = "Perfection"
$ws.Range("I164").Value = "Primera"
$ws.Range("J164").Value = 60
$ws.Range("K164").Value = 32000
$ws.Range("L164").Value = 32000
$ws.Range("M164").Value = 32000
$ws.Range("N164").Value = "$/malla 25 kilos"
$ws.Range("O164").Value = "Provincia de Limarí"
$ws.Range("P164").Value = 1280
$ws.Range("Q164").Value = 25
$ws.Range("R164").Value = "Hortaliza"
